$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the import text block in B2: drop the "import " prefix from the
# first two lines, and change the third line's "import" to "Com".
$ws.Range("B2").Value = "com.blackknight.demo.models.MortgageRequest;`ncom.blackknight.demo.models.Address;`nCom.blackknight.demo.models.Loan;"

# Add new cells B9 and C9 with value "x"
$ws.Range("B9").Value = "x"
$ws.Range("C9").Value = "x"

# Update the active selection to B2
$ws.Range("B2").Select()
